$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Tabelle1" to "Table1"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Table1"

# Add a second sheet "Table2" right after Table1, with the same A1:B3 content
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table2"

$ws2.Range("A1").Value = "product_ID"
$ws2.Range("B1").Value = "username"
$ws2.Range("A2").Value = 253
$ws2.Range("B2").Value = "testUser"
$ws2.Range("A3").Value = 254
$ws2.Range("B3").Value = "testUser2"

# Mirror the author's selection (full data range selected on Table2)
$ws2.Range("A1:B3").Select() | Out-Null

# Move the selection on Table1 back to B19 and re-activate that sheet
$ws1.Range("B19").Select() | Out-Null
$ws1.Activate() | Out-Null
